$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Set-TextPercent($cellRef, $val, $rowNum) {
    # Direct assignment of a pure "NN%" string gets auto-parsed as a percentage number.
    # Force text by prefixing with an apostrophe (quote-prefix), then restore the
    # original (General, s=3) cell style by pasting formats from the same-row
    # URL_FONT cell (column F), which always shares that same style.
    $ws.Range($cellRef).Value = "'" + $val
    $fmtSource = $ws.Range("F" + $rowNum)
    $fmtSource.Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

Set-TextValue "E2" "2026-02-25 19:18:35"
Set-TextValue "O2" "5.7 °C"
Set-TextValue "E3" "2026-02-25 19:18:37"
Set-TextPercent "H3" "36%" 3
Set-TextValue "E4" "2026-02-25 19:18:40"
Set-TextValue "J4" "1021.7 hPa"
Set-TextValue "O4" "8.4 °C"
Set-TextValue "E5" "2026-02-25 19:18:43"
Set-TextValue "N5" "2.8 °C 18:43 TU"
Set-TextValue "O5" "6.0 °C"
Set-TextValue "E6" "2026-02-25 19:18:45"
Set-TextValue "J6" "1021.5 hPa"
Set-TextValue "K6" "4.6 MJ/m2"
Set-TextValue "E7" "2026-02-25 19:18:48"
Set-TextValue "J7" "1021.1 hPa"
Set-TextValue "E8" "2026-02-25 19:18:51"
Set-TextPercent "H8" "77%" 8
Set-TextValue "J8" "1020.7 hPa"
Set-TextValue "E9" "2026-02-25 19:18:53"
Set-TextValue "E10" "2026-02-25 19:18:56"
Set-TextPercent "H10" "89%" 10
Set-TextValue "E11" "2026-02-25 19:18:58"
Set-TextValue "E12" "2026-02-25 19:19:01"
Set-TextValue "E13" "2026-02-25 19:19:03"
Set-TextValue "O13" "6.7 °C"
Set-TextValue "E14" "2026-02-25 19:19:06"
Set-TextPercent "H14" "90%" 14
Set-TextValue "O14" "10.5 °C"
Set-TextValue "E15" "2026-02-25 19:19:09"
Set-TextPercent "H15" "87%" 15
Set-TextValue "O15" "9.9 °C"
Set-TextValue "E16" "2026-02-25 19:19:11"
Set-TextValue "N16" "1.4 °C 18:32 TU"
Set-TextValue "O16" "3.3 °C"
Set-TextValue "E17" "2026-02-25 19:19:14"
Set-TextValue "G17" "1 cm"
Set-TextValue "N17" "6.0 °C 18:54 TU"
Set-TextValue "O17" "9.5 °C"
Set-TextValue "E18" "2026-02-25 19:19:17"
Set-TextValue "O18" "9.9 °C"
Set-TextValue "E19" "2026-02-25 19:19:19"
Set-TextValue "O19" "12.6 °C"
Set-TextValue "E20" "2026-02-25 19:19:22"
Set-TextValue "E21" "2026-02-25 19:19:24"
Set-TextPercent "H21" "56%" 21
Set-TextValue "J21" "1021.4 hPa"
Set-TextValue "E22" "2026-02-25 19:19:27"
Set-TextPercent "H22" "42%" 22
Set-TextValue "E23" "2026-02-25 19:19:29"
Set-TextPercent "H23" "32%" 23
Set-TextValue "O23" "4.1 °C"
Set-TextValue "E24" "2026-02-25 19:19:32"
Set-TextValue "E25" "2026-02-25 19:19:35"
Set-TextValue "K25" "17.3 MJ/m2"
Set-TextValue "E26" "2026-02-25 19:19:37"
Set-TextPercent "H26" "47%" 26
Set-TextValue "J26" "1019.1 hPa"
Set-TextValue "N26" "6.5 °C 18:59 TU"
Set-TextValue "O26" "10.5 °C"
Set-TextValue "E27" "2026-02-25 19:19:40"
Set-TextValue "K27" "16.5 MJ/m2"
Set-TextValue "E28" "2026-02-25 19:19:43"
Set-TextPercent "H28" "84%" 28
Set-TextValue "J28" "1021.6 hPa"
Set-TextValue "O28" "8.7 °C"
Set-TextValue "E29" "2026-02-25 19:19:46"
Set-TextValue "O29" "12.0 °C"
Set-TextValue "E30" "2026-02-25 19:19:48"
Set-TextValue "J30" "1021.6 hPa"
Set-TextValue "O30" "10.5 °C"
Set-TextValue "E31" "2026-02-25 19:19:51"
Set-TextValue "J31" "1021.2 hPa"
Set-TextValue "E32" "2026-02-25 19:19:54"
Set-TextValue "O32" "9.8 °C"
Set-TextValue "E33" "2026-02-25 19:19:56"
Set-TextValue "J33" "1021.0 hPa"
Set-TextValue "E34" "2026-02-25 19:19:59"
Set-TextPercent "H34" "51%" 34
Set-TextValue "E35" "2026-02-25 19:20:01"
Set-TextValue "J35" "1019.2 hPa"
Set-TextValue "O35" "12.7 °C"
Set-TextValue "E36" "2026-02-25 19:20:04"
Set-TextValue "J36" "1021.7 hPa"
Set-TextValue "O36" "11.1 °C"
Set-TextValue "E37" "2026-02-25 19:20:06"
Set-TextPercent "H37" "82%" 37
Set-TextValue "J37" "1023.1 hPa"
Set-TextValue "O37" "6.9 °C"
Set-TextValue "E38" "2026-02-25 19:20:09"
Set-TextValue "E39" "2026-02-25 19:20:11"
Set-TextPercent "H39" "45%" 39
Set-TextValue "E40" "2026-02-25 19:20:14"
Set-TextValue "E41" "2026-02-25 19:20:16"
Set-TextValue "J41" "1020.7 hPa"
Set-TextValue "E42" "2026-02-25 19:20:19"
Set-TextValue "E43" "2026-02-25 19:20:21"
Set-TextValue "O43" "10.0 °C"
Set-TextValue "E44" "2026-02-25 19:20:24"
Set-TextValue "E45" "2026-02-25 19:20:27"
Set-TextValue "L45" "21.6 km/h - 110º 18:33 TU"
Set-TextValue "E46" "2026-02-25 19:20:29"
Set-TextValue "O46" "9.6 °C"
